# Updates results table "table_parameters_irt" for revisions with updated Stan version.
# Only the numeric mean/sd/q05/q50/q95 values change; row/column labels are untouched.
# Cell values are written with a leading apostrophe so Excel stores them as text
# (matching the original shared-string/text cells, incl. preserved leading spaces),
# then the style is reset to "Normal" so no stray text-format style is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ($\beta_1$)
$ws.Range("B2").Value = "'-0.766"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "'1.009"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'-2.522"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.673"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'0.698"
$ws.Range("F2").Style = "Normal"
# Row 3 ($\beta_2$)
$ws.Range("B3").Value = "' 1.660"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'0.926"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "' 0.183"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "' 1.631"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'3.205"
$ws.Range("F3").Style = "Normal"
# Row 4 ($\alpha$)
$ws.Range("B4").Value = "' 0.112"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'0.320"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'-0.417"
$ws.Range("D4").Style = "Normal"
$ws.Range("F4").Value = "'0.616"
$ws.Range("F4").Style = "Normal"
# Row 5 ($\mu_1$)
$ws.Range("B5").Value = "' 0.557"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'0.975"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "'-1.069"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "' 0.561"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'2.146"
$ws.Range("F5").Style = "Normal"
# Row 6 ($\mu_2$)
$ws.Range("B6").Value = "' 2.557"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'1.007"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "' 0.971"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "' 2.534"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'4.237"
$ws.Range("F6").Style = "Normal"
# Row 7 ($\mu_3$)
$ws.Range("B7").Value = "' 1.357"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'1.034"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'-0.309"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "' 1.337"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'3.087"
$ws.Range("F7").Style = "Normal"
# Row 8 ($\mu_4$)
$ws.Range("B8").Value = "' 2.253"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'1.086"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "' 0.490"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "' 2.235"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'4.061"
$ws.Range("F8").Style = "Normal"
# Row 9 ($\sigma_b$)
$ws.Range("B9").Value = "' 2.200"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'0.663"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "' 1.308"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "' 2.106"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'3.436"
$ws.Range("F9").Style = "Normal"
# Row 10 ($\sigma_a$)
$ws.Range("B10").Value = "' 0.219"
$ws.Range("B10").Style = "Normal"
$ws.Range("D10").Value = "' 0.011"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "' 0.129"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'0.835"
$ws.Range("F10").Style = "Normal"
# Row 13 ($\rho_{a,b}$)
$ws.Range("B13").Value = "'-0.122"
$ws.Range("B13").Style = "Normal"
$ws.Range("D13").Value = "'-0.899"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.182"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = "'0.816"
$ws.Range("F13").Style = "Normal"
# Row 14 ($\rho_{Gov,Out}$)
$ws.Range("C14").Value = "'0.346"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'-0.139"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "' 0.594"
$ws.Range("E14").Style = "Normal"

Write-Host "table_parameters_irt values updated"
